$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 302.3
$ws.Range("H40").Value = 6760.1333
$ws.Range("I40").Value = 4640.3
$ws.Range("J40").Value = 10999.8
$ws.Range("K40").Value = 4640.3
$ws.Range("L40").Value = 10999.8
$ws.Range("M40").Value = -4465.3
$ws.Range("N40").Value = -11349.8
$ws.Range("H107").Value = 78712.54
$ws.Range("I107").Value = 78712.54
$ws.Range("K107").Value = 78712.54
$ws.Range("M107").Value = -76792.54
$ws.Range("H112").Value = 2854.182
$ws.Range("I112").Value = 100
$ws.Range("J112").Value = 2985.3333
$ws.Range("K112").Value = 300
$ws.Range("L112").Value = 8955.999899999999
$ws.Range("M112").Value = 808
$ws.Range("N112").Value = -11171.9999
$ws.Range("H125").Value = 12348063
$ws.Range("I125").Value = 791.2
$ws.Range("J125").Value = 27782152
$ws.Range("K125").Value = 7120.8
$ws.Range("L125").Value = 250039368
$ws.Range("M125").Value = -4660.8
$ws.Range("N125").Value = -250044288
$ws.Range("H132").Value = 2620.5
$ws.Range("I132").Value = 2432
$ws.Range("K132").Value = 7296
$ws.Range("M132").Value = -4766
$ws.Range("H135").Value = 716753.4
$ws.Range("J135").Value = 2052
$ws.Range("L135").Value = 18468
$ws.Range("N135").Value = -23538
$ws.Range("H138").Value = 5233.174
$ws.Range("I138").Value = 2784.7144
$ws.Range("J138").Value = 5672.641
$ws.Range("K138").Value = 8354.143199999999
$ws.Range("L138").Value = 17017.923
$ws.Range("M138").Value = -3214.143199999999
$ws.Range("N138").Value = -27297.923

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1492.44
$ws.Range("I32").Value = 1492.44
$ws.Range("K32").Value = 1492.44
$ws.Range("M32").Value = -1205.44
$ws.Range("H45").Value = 2027
$ws.Range("I45").Value = 1569.3334
$ws.Range("K45").Value = 1569.3334
$ws.Range("M45").Value = -1192.3334
$ws.Range("H61").Value = 2785.5715
$ws.Range("I61").Value = 2846
$ws.Range("K61").Value = 2846
$ws.Range("M61").Value = -2634
$ws.Range("H63").Value = 5988.6924
$ws.Range("I63").Value = 2974.5
$ws.Range("K63").Value = 2974.5
$ws.Range("M63").Value = -2288.5
$ws.Range("H66").Value = 5988.6924
$ws.Range("I66").Value = 2974.5
$ws.Range("K66").Value = 14872.5
$ws.Range("M66").Value = -11440.5
$ws.Range("H132").Value = 2976.762
$ws.Range("I132").Value = 3075.7222
$ws.Range("K132").Value = 9227.1666
$ws.Range("M132").Value = -6697.1666
$ws.Range("H136").Value = 2785.5715
$ws.Range("I136").Value = 2846
$ws.Range("K136").Value = 8538
$ws.Range("M136").Value = -5988
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 22950
$ws.Range("J140").Value = 22950
$ws.Range("L140").Value = 22950
$ws.Range("N140").Value = -33310
$ws.Range("H141").Value = 23950
$ws.Range("J141").Value = 23950
$ws.Range("L141").Value = 23950
$ws.Range("N141").Value = -34310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 15000
$ws.Range("J44").Value = 15000
$ws.Range("L44").Value = 15000
$ws.Range("N44").Value = -15994
$ws.Range("H94").Value = 1566.4117
$ws.Range("I94").Value = 1750.8462
$ws.Range("J94").Value = 967
$ws.Range("K94").Value = 1750.8462
$ws.Range("L94").Value = 967
$ws.Range("M94").Value = -1299.8462
$ws.Range("N94").Value = -1869
$ws.Range("H134").Value = 34385.324
$ws.Range("I134").Value = 4670.033
$ws.Range("J134").Value = 257250
$ws.Range("K134").Value = 14010.099
$ws.Range("L134").Value = 771750
$ws.Range("M134").Value = -11475.099
$ws.Range("N134").Value = -776820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 20002750
$ws.Range("I4").Value = 3001
$ws.Range("J4").Value = 40002500
$ws.Range("K4").Value = 3001
$ws.Range("L4").Value = 40002500
$ws.Range("M4").Value = -2889
$ws.Range("N4").Value = -40002724
$ws.Range("H6").Value = 5125
$ws.Range("I6").Value = 5000
$ws.Range("J6").Value = 6000
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = -4887
$ws.Range("N6").Value = -6226
$ws.Range("H31").Value = 86984.836
$ws.Range("I31").Value = 1851.375
$ws.Range("J31").Value = 257251.75
$ws.Range("K31").Value = 1851.375
$ws.Range("L31").Value = 257251.75
$ws.Range("M31").Value = -1556.375
$ws.Range("N31").Value = -257841.75
$ws.Range("H34").Value = 86984.836
$ws.Range("I34").Value = 1851.375
$ws.Range("J34").Value = 257251.75
$ws.Range("K34").Value = 1851.375
$ws.Range("L34").Value = 257251.75
$ws.Range("M34").Value = -1649.375
$ws.Range("N34").Value = -257655.75
$ws.Range("H59").Value = 31123.666
$ws.Range("J59").Value = 32139.125
$ws.Range("L59").Value = 32139.125
$ws.Range("N59").Value = -34429.125
$ws.Range("H94").Value = 1032.5385
$ws.Range("J94").Value = 1244.4445
$ws.Range("L94").Value = 1244.4445
$ws.Range("N94").Value = -2146.4445
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676
$ws.Range("H132").Value = 3001.25
$ws.Range("I132").Value = 2873.75
$ws.Range("J132").Value = 3256.25
$ws.Range("K132").Value = 8621.25
$ws.Range("L132").Value = 9768.75
$ws.Range("M132").Value = -6091.25
$ws.Range("N132").Value = -14828.75
$ws.Range("H134").Value = 387305.84
$ws.Range("I134").Value = 2755.7827
$ws.Range("K134").Value = 8267.348100000001
$ws.Range("M134").Value = -5732.348100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 128584.375
$ws.Range("J34").Value = 146925
$ws.Range("L34").Value = 440775
$ws.Range("N34").Value = -440943
$ws.Range("H80").Value = 1200.6
$ws.Range("I80").Value = 1501
$ws.Range("J80").Value = 1000.3333
$ws.Range("K80").Value = 4503
$ws.Range("L80").Value = 3000.9999
$ws.Range("M80").Value = -3567
$ws.Range("N80").Value = -4872.9999
$ws.Range("H83").Value = 1200.6
$ws.Range("I83").Value = 1501
$ws.Range("J83").Value = 1000.3333
$ws.Range("K83").Value = 13509
$ws.Range("L83").Value = 9002.9997
$ws.Range("M83").Value = -8829
$ws.Range("N83").Value = -18362.9997
$ws.Range("H92").Value = 909925.06
$ws.Range("J92").Value = 2750
$ws.Range("L92").Value = 8250
$ws.Range("N92").Value = -10746
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H98").Value = 4236.1665
$ws.Range("J98").Value = 3621.75
$ws.Range("L98").Value = 10865.25
$ws.Range("N98").Value = -13861.25
$ws.Range("H101").Value = 6352.6665
$ws.Range("J101").Value = 7779
$ws.Range("L101").Value = 23337
$ws.Range("N101").Value = -28205
$ws.Range("H109").Value = 48842.316
$ws.Range("I109").Value = 1636.7778
$ws.Range("J109").Value = 81523.08
$ws.Range("K109").Value = 4910.3334
$ws.Range("L109").Value = 244569.24
$ws.Range("M109").Value = -3870.3334
$ws.Range("N109").Value = -246649.24
$ws.Range("H140").Value = 4637
$ws.Range("I140").Value = 4129.5
$ws.Range("K140").Value = 12388.5
$ws.Range("M140").Value = -7208.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3942.1428
$ws.Range("I122").Value = 1398.75
$ws.Range("K122").Value = 4196.25
$ws.Range("M122").Value = -1746.25
$ws.Range("H132").Value = 99144.91
$ws.Range("I132").Value = 9510.444
$ws.Range("K132").Value = 28531.332
$ws.Range("M132").Value = -26001.332
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 59161.332
$ws.Range("J139").Value = 59161.332
$ws.Range("L139").Value = 59161.332
$ws.Range("N139").Value = -69441.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3795.75
$ws.Range("I61").Value = 3795.75
$ws.Range("K61").Value = 3795.75
$ws.Range("M61").Value = -3593.75
$ws.Range("H82").Value = 4746.3335
$ws.Range("J82").Value = 4747
$ws.Range("L82").Value = 4747
$ws.Range("N82").Value = -5469
$ws.Range("H85").Value = 4746.3335
$ws.Range("J85").Value = 4747
$ws.Range("L85").Value = 4747
$ws.Range("N85").Value = -7243
$ws.Range("H100").Value = 3700.6667
$ws.Range("J100").Value = 3002
$ws.Range("L100").Value = 3002
$ws.Range("N100").Value = -4084
$ws.Range("H113").Value = 3795.75
$ws.Range("I113").Value = 3795.75
$ws.Range("K113").Value = 3795.75
$ws.Range("M113").Value = -1625.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H122").Value = 32261796
$ws.Range("I122").Value = 37040492
$ws.Range("J122").Value = 5601.25
$ws.Range("K122").Value = 111121476
$ws.Range("L122").Value = 16803.75
$ws.Range("M122").Value = -111119026
$ws.Range("N122").Value = -21703.75
$ws.Range("H132").Value = 44075.6
$ws.Range("I132").Value = 3136.0527
$ws.Range("J132").Value = 173717.5
$ws.Range("K132").Value = 9408.1581
$ws.Range("L132").Value = 521152.5
$ws.Range("M132").Value = -6878.158100000001
$ws.Range("N132").Value = -526212.5
$ws.Range("H136").Value = 222284.1
$ws.Range("I136").Value = 44512.5
$ws.Range("K136").Value = 133537.5
$ws.Range("M136").Value = -130987.5
